# Allow not to press Enter when finish typing the name option in ketxuat
#
# Adds a new "ton-thanhpham" (stock-on-hand) worksheet, mirroring the
# nhap-thanhpham / xuat-thanhpham sheets (same headers + column widths),
# and seeds it with a single summary row: item "he" with a running total
# of -100 in the Số Lượng (quantity) column.

$wb = $excel.ActiveWorkbook

# Copy an existing sheet (right after itself) so the header row, column
# widths and page setup all match the other two sheets, then rename it.
$srcSheet = $wb.Worksheets.Item("xuat-thanhpham")
$srcSheet.Copy($null, $srcSheet)
$newSheet = $wb.ActiveSheet
$newSheet.Name = "ton-thanhpham"

# Drop the copied data row and replace it with the ton-thanhpham summary:
# only the item name and the (negative) quantity are populated.
$newSheet.Rows.Item(2).ClearContents()
$newSheet.Range("A2").Value = "he"
$newSheet.Range("F2").Value = -100
